$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has data in rows 1,3,4,7-13 with blank rows at 2,5,6.
# Remove the blank rows (2, 5, 6) so everything shifts up and becomes a
# contiguous block in rows 1-10. Delete from the bottom up so row numbers
# of not-yet-deleted rows stay valid.
$ws.Rows("5:6").Delete()
$ws.Rows("2:2").Delete()

# Dimension is now A1:B10. Update the selection to match the new used range
# (whole block selected, no single active-cell override).
$ws.Range("A1:B10").Select()
